$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-51 price/volume updates from the Nov 27 2023 GitHub Actions crypto data refresh.
# Column D (Price) cells are stored as text in this sheet, so values are entered with a
# leading apostrophe to force text (matching the existing t="inlineStr" cell type) instead
# of letting Excel auto-convert numeric-looking strings into real numbers.

$ws.Range("D2").Value = "'36.827.36"
$ws.Range("E2").Value = "  -1.72%  "

$ws.Range("D3").Value = "'2.018.98"
$ws.Range("E3").Value = "  -2.74%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'225.69"
$ws.Range("E5").Value = "  -2.84%  "

$ws.Range("D6").Value = "'0.605"
$ws.Range("E6").Value = "  -4.04%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'54.60"
$ws.Range("E8").Value = "  -5.10%  "

$ws.Range("D9").Value = "'0.379"
$ws.Range("E9").Value = "  -2.44%  "

$ws.Range("D10").Value = "'0.0786"
$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("E11").Value = "  -3.69%  "

$ws.Range("D12").Value = "'2.319.53"
$ws.Range("E12").Value = "  -2.73%  "

$ws.Range("D13").Value = "'14.24"
$ws.Range("E13").Value = "  -4.45%  "

$ws.Range("D14").Value = "'20.32"
$ws.Range("E14").Value = "  -2.99%  "

$ws.Range("D15").Value = "'0.744"
$ws.Range("E15").Value = "  -3.13%  "

$ws.Range("D16").Value = "'5.13"
$ws.Range("E16").Value = "  -3.85%  "

$ws.Range("D17").Value = "'2.020.79"
$ws.Range("E17").Value = "  -2.45%  "

$ws.Range("D18").Value = "'36.827.68"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("D19").Value = "'6.17"
$ws.Range("E19").Value = "  +2.80%  "

$ws.Range("D20").Value = "'68.60"
$ws.Range("E20").Value = "  -2.65%  "

$ws.Range("D21").Value = "'0.0₃0820"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("D22").Value = "'225.49"
$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  +3.07%  "

$ws.Range("D25").Value = "'2.24"
$ws.Range("E25").Value = "  -5.24%  "

$ws.Range("D26").Value = "'165.54"
$ws.Range("E26").Value = "  -2.62%  "

$ws.Range("D27").Value = "'9.16"
$ws.Range("E27").Value = "  -4.88%  "

$ws.Range("D28").Value = "'0.124"
$ws.Range("E28").Value = "  -5.85%  "

$ws.Range("D29").Value = "'18.67"
$ws.Range("E29").Value = "  -4.27%  "

$ws.Range("D30").Value = "'1.32"
$ws.Range("E30").Value = "  -3.43%  "

$ws.Range("E31").Value = "  -5.16%  "

$ws.Range("D32").Value = "'4.45"
$ws.Range("E32").Value = "  -3.60%  "

$ws.Range("D33").Value = "'0.0614"
$ws.Range("E33").Value = "  -2.74%  "

$ws.Range("D34").Value = "'4.41"
$ws.Range("E34").Value = "  -4.66%  "

$ws.Range("D35").Value = "'2.35"
$ws.Range("E35").Value = "  -4.80%  "

$ws.Range("D36").Value = "'1.83"
$ws.Range("E36").Value = "  +1.15%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").Value = "'3.14"
$ws.Range("E38").Value = "  -5.02%  "

$ws.Range("D39").Value = "'5.27"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'1.488.32"
$ws.Range("E40").Value = "  +1.92%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0218"
$ws.Range("E41").Value = "  -5.23%  "

$ws.Range("D42").Value = "'16.95"
$ws.Range("E42").Value = "  +1.92%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0926"
$ws.Range("E43").Value = "  -2.91%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'94.73"
$ws.Range("E44").Value = "  -5.16%  "

$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "'2.74"
$ws.Range("E45").Value = "  -6.11%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.13"
$ws.Range("E46").Value = "  -5.79%  "

$ws.Range("D47").Value = "'7.30"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -4.11%  "

$ws.Range("D49").Value = "'2.90"
$ws.Range("E49").Value = "  -1.54%  "

$ws.Range("D50").Value = "'2.209.85"
$ws.Range("E50").Value = "  -2.58%  "

$ws.Range("D51").Value = "'3.58"
$ws.Range("E51").Value = "  -9.28%  "
